$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BCoESC-power-plants")

$newPlants = @(
    "hard coal w CCS",
    "natural gas combined cycle w CCS",
    "biomass w CCS",
    "lignite w CCS",
    "small modular reactor",
    "hydrogen"
)

$row = 19
foreach ($plant in $newPlants) {
    $ws.Cells.Item($row, 1).Value = $plant
    $ws.Cells.Item($row, 2).Formula = '=Data!$B$20'
    $row = $row + 1
}

$ws.Range("A25").Select()
